$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B3").Value = 9
$ws.Range("B6").Value = 109
